$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataA = @(5310,5260,5230,5200,5180,5160,5140,5120,5110,5100,5100,5120,5130,5140,5160,5170,5190,5220,5260,5330,5380,5480,5590,5720,5900,6030,6120,6140,6140,6140,6140,6100,6030,5950,5860,5770,5690,5610,5550,5470,5390,5340,5290,5260,5220,5190,5170,5160,5190,5200,5200,5200,5200,5200,5200,5230,5260,5280,5310,5350,5400,5480,5550,5630,5720,5810,5890,5970,6090,6170,6260,6360,6440,6550,6660,6760,6900,7010,7090,7160,7180,7120,7060,6950,6770,6590,6420,6240,6090,5910,5810,5700,5530,5490,5430,5350)
$dataB = @(45786,45786.01041666666,45786.02083333334,45786.03125,45786.04166666666,45786.05208333334,45786.0625,45786.07291666666,45786.08333333334,45786.09375,45786.10416666666,45786.11458333334,45786.125,45786.13541666666,45786.14583333334,45786.15625,45786.16666666666,45786.17708333334,45786.1875,45786.19791666666,45786.20833333334,45786.21875,45786.22916666666,45786.23958333334,45786.25,45786.26041666666,45786.27083333334,45786.28125,45786.29166666666,45786.30208333334,45786.3125,45786.32291666666,45786.33333333334,45786.34375,45786.35416666666,45786.36458333334,45786.375,45786.38541666666,45786.39583333334,45786.40625,45786.41666666666,45786.42708333334,45786.4375,45786.44791666666,45786.45833333334,45786.46875,45786.47916666666,45786.48958333334,45786.5,45786.51041666666,45786.52083333334,45786.53125,45786.54166666666,45786.55208333334,45786.5625,45786.57291666666,45786.58333333334,45786.59375,45786.60416666666,45786.61458333334,45786.625,45786.63541666666,45786.64583333334,45786.65625,45786.66666666666,45786.67708333334,45786.6875,45786.69791666666,45786.70833333334,45786.71875,45786.72916666666,45786.73958333334,45786.75,45786.76041666666,45786.77083333334,45786.78125,45786.79166666666,45786.80208333334,45786.8125,45786.82291666666,45786.83333333334,45786.84375,45786.85416666666,45786.86458333334,45786.875,45786.88541666666,45786.89583333334,45786.90625,45786.91666666666,45786.92708333334,45786.9375,45786.94791666666,45786.95833333334,45786.96875,45786.97916666666,45786.98958333334)

for ($i = 0; $i -lt $dataA.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dataA[$i]
    $ws.Cells.Item($row, 2).Value = $dataB[$i]
}
